$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Gen slack" sheet: update column B values (rows 2-4) and selection
# ---------------------------------------------------------------------------
$wsGen = $wb.Worksheets.Item("Gen slack")
$wsGen.Range("B2").Value = 1
$wsGen.Range("B3").Value = 4
$wsGen.Range("B4").Value = 9

# ---------------------------------------------------------------------------
# 2. "Bus" sheet: update column B values (rows 4-10) and append new rows
#    11-14 (with column A carrying the same centred style as A10).
# ---------------------------------------------------------------------------
$wsBus = $wb.Worksheets.Item("Bus")
$wsBus.Range("B4").Value = 11
$wsBus.Range("B5").Value = 225
$wsBus.Range("B6").Value = 33
$wsBus.Range("B7").Value = 33
$wsBus.Range("B8").Value = 11
$wsBus.Range("B10").Value = 225

$wsBus.Range("A10").Copy($wsBus.Range("A11"))
$wsBus.Range("A10").Copy($wsBus.Range("A12"))
$wsBus.Range("A10").Copy($wsBus.Range("A13"))
$wsBus.Range("A10").Copy($wsBus.Range("A14"))

$wsBus.Range("A11").Value = 9
$wsBus.Range("B11").Value = 30
$wsBus.Range("A12").Value = 10
$wsBus.Range("B12").Value = 30
$wsBus.Range("A13").Value = 11
$wsBus.Range("B13").Value = 11
$wsBus.Range("A14").Value = 12
$wsBus.Range("B14").Value = 11

# ---------------------------------------------------------------------------
# 3. "Lines" sheet: rebuild the data block so that it holds 10 rows (2-11)
#    instead of 5, dropping the stray row 30. The two rows that originally
#    carried the leftover N:U formatting (rows 5 & 6) are pushed down to
#    rows 6 & 7 by inserting a row above the old row 5, then four more rows
#    are appended right after for the new entries.
# ---------------------------------------------------------------------------
$wsLines = $wb.Worksheets.Item("Lines")
$wsLines.Rows.Item(30).Delete()
$wsLines.Range("5:5").Insert()
$wsLines.Range("8:11").Insert()
$wsLines.Range("N8:U11").Clear()

$linesData = @(
  @(2, 0, 1, 2, 0, 0, 0, 0, 0),
  @(3, 1, 2, 3, 125, 0.6, 0.08, 210, 0.142),
  @(4, 2, 4, 5, 0, 0, 0, 0, 0),
  @(5, 3, 5, 6, 65, 0.6, 0.08, 210, 0.142),
  @(6, 4, 6, 7, 0, 0, 0, 0, 0),
  @(7, 5, 7, 8, 60, 0.6, 0.08, 210, 0.142),
  @(8, 6, 9, 10, 0, 0, 0, 0, 0),
  @(9, 7, 10, 11, 65, 0.6, 0.08, 210, 0.142),
  @(10, 8, 11, 12, 0, 0, 0, 0, 0),
  @(11, 9, 12, 13, 60, 0.6, 0.08, 210, 0.142)
)
foreach ($r in $linesData) {
  $rn = $r[0]
  $wsLines.Cells.Item($rn, 1).Value = $r[1]
  $wsLines.Cells.Item($rn, 2).Value = $r[2]
  $wsLines.Cells.Item($rn, 3).Value = $r[3]
  $wsLines.Cells.Item($rn, 4).Value = $r[4]
  $wsLines.Cells.Item($rn, 5).Value = $r[5]
  $wsLines.Cells.Item($rn, 6).Value = $r[6]
  $wsLines.Cells.Item($rn, 7).Value = $r[7]
  $wsLines.Cells.Item($rn, 8).Value = $r[8]
}

# ---------------------------------------------------------------------------
# 4. Remove the "line1" and "Line2" sheets entirely.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("line1").Delete()
$wb.Worksheets.Item("Line2").Delete()

# ---------------------------------------------------------------------------
# 5. "Load" sheet (the former last sheet): update/append rows.
# ---------------------------------------------------------------------------
$wsLoad = $wb.Worksheets.Item("Load")
$wsLoad.Range("B2").Value = 3
$wsLoad.Range("A3").Value = 1
$wsLoad.Range("B3").Value = 8
$wsLoad.Range("C3").Value = 0.03
$wsLoad.Range("A4").Value = 2
$wsLoad.Range("B4").Value = 13
$wsLoad.Range("C4").Value = 0.03

# ---------------------------------------------------------------------------
# 6. Selections / active sheet bookkeeping (must run after the sheet
#    deletions above so the tab order & indices line up with the target).
# ---------------------------------------------------------------------------
$wsGen.Range("F11").Select()
$wsBus.Range("H15").Select()
$wsLines.Range("K25").Select()

$wsLoad.Activate()
$wsLoad.Range("H14").Select()
